$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, matching the header style used by the existing header row (A1:AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in team record data for each data row (2-37)
$ws.Range("AD2:AD37").Value = 116
$ws.Range("AE2:AE37").Value = 46
$ws.Range("AF2:AF37").Value = 0
